$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val -is [string] -and $val.EndsWith(".jpg")) {
            $newVal = $val.Substring(0, $val.Length - 4) + ".png"
            $cell.Value2 = $newVal
        }
    }
}

# Reset the selection to the full used range and clear the frozen/scroll position
[void]$used.Select()
